$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.061.57'
$ws.Range("E2").Value = '  +0.67%  '

$ws.Range("D3").Value = '2.352.63'
$ws.Range("E3").Value = '  +0.15%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.678'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '238.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.63%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.91'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.75%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.591'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +9.09%  '

$ws.Range("E10").Value = '  +2.11%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.19'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.27%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '32.14'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +13.68%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.108'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.87%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.14'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.71%  '

$ws.Range("D15").Value = '2.697.86'
$ws.Range("E15").Value = '  -0.09%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '16.51'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.00%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.898'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.26%  '

$ws.Range("D18").Value = '2.358.14'
$ws.Range("E18").Value = '  -0.13%  '

$ws.Range("D19").Value = '43.905.55'
$ws.Range("E19").Value = '  +0.42%  '

$ws.Range("E20").Value = '  +1.11%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.67'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.50%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '76.68'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.32%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '256.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.36%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +20.24%  '

$ws.Range("E25").Value = '  +0.01%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.69'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.59%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.48'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.24%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.65'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.89%  '

$ws.Range("E29").Value = '  -2.36%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.75'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.88%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.18'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.44%  '

$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.137'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.70%  '

$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.126'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.64%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0757'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.20%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.27'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.32'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.62%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.72'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.52%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.36'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.11%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.32'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.83%  '

$ws.Range("E40").Value = '  +4.52%  '

$ws.Range("E41").Value = '  +11.83%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '19.11'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.44%  '

$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.202'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +10.72%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.03'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.89%  '

$ws.Range("E45").Value = '  -0.05%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.68'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.62%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.48'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +8.96%  '

$ws.Range("E48").Value = '  +1.71%  '

$ws.Range("B49").Value = 'MultiversX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '57.03'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.21%  '

$ws.Range("E50").Value = '  +1.15%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '99.89'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.92%  '
